# DS-2197: Rename attribute level "Belgian (Single origin Venezuelan
# Criollo beans)" -> "Belgium (Single origin Venezuelan Criollo beans)"
# on the "attributes" sheet, and leave the selection on C5 (matching the
# state Excel was left in when the workbook was saved).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value2 = "Belgium (Single origin Venezuelan Criollo beans)"

$ws.Range("C5").Select()
